$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# --- Paragraph 8 (bookmark paragraph) -> add rPr lang=en-US, then append a
# new trailing empty paragraph (ListParagraph, ind left=1440, no numPr) ---
$p8 = $d.Paragraphs.Item(8)
$frag8 = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"1`"/><w:numId w:val=`"1`"/></w:numPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/></w:p><w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:ind w:left=`"1440`"/></w:pPr></w:p>"
$p8.Range.InsertXML($frag8)

# --- Paragraphs 5-7 (Promedio diario.../ghi_diario_mensual/ghi_diario) ->
# replaced by four new "Valledupar 2010 - N mes" bullets ---
$p5 = $d.Paragraphs.Item(5)
$p7 = $d.Paragraphs.Item(7)
$rng567 = $d.Range($p5.Range.Start, $p7.Range.End)

$fragFirst = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"1`"/><w:numId w:val=`"1`"/></w:numPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`">Valledupar 2010 – 1 </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>mes</w:t></w:r><w:proofErr w:type=`"spellEnd`"/></w:p>"

$frag3m = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"1`"/><w:numId w:val=`"1`"/></w:numPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`">Valledupar 2010 – </w:t></w:r><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>3</w:t></w:r><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>mes</w:t></w:r><w:proofErr w:type=`"spellEnd`"/></w:p>"

$frag6m = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"1`"/><w:numId w:val=`"1`"/></w:numPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`">Valledupar 2010 – </w:t></w:r><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>6</w:t></w:r><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>mes</w:t></w:r><w:proofErr w:type=`"spellEnd`"/></w:p>"

$frag9m = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"1`"/><w:numId w:val=`"1`"/></w:numPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`">Valledupar 2010 – </w:t></w:r><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>9</w:t></w:r><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>mes</w:t></w:r><w:proofErr w:type=`"spellEnd`"/></w:p>"

$frag567 = $fragFirst + $frag3m + $frag6m + $frag9m
$rng567.InsertXML($frag567)

# --- Paragraph 4 (MCP: mcp) -> "Valledupar 2010" ---
$p4 = $d.Paragraphs.Item(4)
$frag4 = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"1`"/><w:numId w:val=`"1`"/></w:numPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>Valledupar 2010</w:t></w:r></w:p>"
$p4.Range.InsertXML($frag4)

# --- Paragraph 3 (Cortar series: cut_series) -> "Informe" at ilvl 0 ---
$p3 = $d.Paragraphs.Item(3)
$frag3 = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>Informe</w:t></w:r></w:p>"
$p3.Range.InsertXML($frag3)
